$wb = $excel.ActiveWorkbook

# --- Sheet "ZFCode" (3rd sheet) ---
$ws3 = $wb.Worksheets.Item(3)

# Highlight B1:B3 and F1:F3 with a yellow fill
$ws3.Range("B1:B3").Interior.Color = 65535
$ws3.Range("F1:F3").Interior.Color = 65535

# Note about retrieving from PMT daily data
$ws3.Cells.Item(2, 6).Value = "retrieve from PMT daily data"

# New entry at the bottom referencing the sales rev turn retrieval
$ws3.Cells.Item(18, 2).Value = "retrieve_sales_rev_turn_raw"

$ws3.PageSetup.Orientation = 1

# --- Sheet "Julia tip" (2nd sheet) ---
$ws2 = $wb.Worksheets.Item(2)

$ws2.Cells.Item(21, 2).Value = "cat multiple variables"
$ws2.Cells.Item(21, 3).Value = "cat(x…)"

$ws2.Cells.Item(22, 1).Value = "replace."
$ws2.Cells.Item(22, 2).Value = "cannot broadcast"
$ws2.Cells.Item(22, 3).Value = "map is a quite good tool"

$ws2.Cells.Item(24, 1).Value = "LinearIndices"
$ws2.Cells.Item(24, 2).Value = "CartesianIndex"
$ws2.Cells.Item(24, 3).Value = "transfer"

$ws2.Cells.Item(26, 1).Value = "Debug"
$rngB26 = $ws2.Cells.Item(26, 2)
$rngB26.Value = "Juno.@enter function(inputs)"
$rngB26.Font.Name = "Inherit"
$rngB26.Font.Color = 0
$rngB26.HorizontalAlignment = -4131
$rngB26.VerticalAlignment = -4108

# --- Sheet "Changes from matlab" (1st sheet) ---
$ws1 = $wb.Worksheets.Item(1)

$ws1.Cells.Item(10, 2).Value = "get_country_sizeInfo"
$ws1.Cells.Item(10, 4).Value = "Sales_rev_turn_raw mat file to jld file"
$ws1.Cells.Item(11, 4).Value = "salesRevTurnMth mat file to jld file"

# --- Selections / active sheet ---
$ws3.Range("D27").Select() | Out-Null
$ws2.Range("C25").Select() | Out-Null

$ws1.Activate() | Out-Null
$ws1.Range("C17:C18").Select() | Out-Null
